# DDAf_2023_Tableau_annexe_Tab10 - "Add files via upload" data refresh.
# The underlying statistics for the "fragile states" aggregate rows (and one
# upstream regional-average cell they feed into) were recomputed; this
# script pokes the refreshed figures into the Tab10 worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab10")

# Regional average row (Afrique, pays à revenu intermediaire tranche
# supérieure) - column I picked up a tiny downstream correction.
$ws.Range("I68").Value = 45.645454545454598

# Row 97: "Afrique, États fragiles" - refreshed aggregate.
$ws.Range("C97").Value = 0.52607142857142997
$ws.Range("D97").Value = 0.88707692307691999
$ws.Range("E97").Value = 0.56604545454544997
$ws.Range("F97").Value = 39.8291666666667
$ws.Range("G97").Value = 52.8392857142857
$ws.Range("H97").Value = 31.02
$ws.Range("I97").Value = 41.588
$ws.Range("J97").Value = 34.5703703703704

# Row 98: "RDM, États fragiles" - refreshed aggregate.
$ws.Range("C98").Value = 0.63606666666667
$ws.Range("D98").Value = 0.91858333333332998
$ws.Range("E98").Value = 0.48275000000000001
$ws.Range("F98").Value = 36.0363636363636
$ws.Range("G98").Value = 41.4
$ws.Range("H98").Value = 34.6133333333333
$ws.Range("I98").Value = 35.2
$ws.Range("J98").Value = 29.0727272727273
